$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "Solution1:" + [char]10 + "1. Ask client for microsoft error debugging."
$ws.Range("C3").Value = "Solution2:" + [char]10 + "1. Ask client for Security token ."

$ws.Range("C3").Select()
